$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that no longer belong in the list (operadoras removed).
# Deleting from bottom to top so earlier row numbers stay valid.
$ws.Rows(63).Delete()   # WIZEO
$ws.Rows(43).Delete()   # RAPPI FARMACIA
$ws.Rows(6).Delete()    # BENVISAVALE
